$d = $word.ActiveDocument

# The edit reorders several paragraphs of the "Objetivos" / "Docente(s)
# Responsavel(eis)" / "Programa resumido" / "Programa" / "Avaliacao"
# sections, changes a couple of paragraph styles along the way, and
# rotates text between runs inside the "Avaliacao" and "Bibliografia"
# paragraphs. Rather than chase that through Find/Replace (several of
# the strings recur / shift positions, which makes ordinary
# find-and-replace ambiguous), each affected paragraph is rewritten in
# place via Range.InsertXML with the exact WordprocessingML it should
# contain afterwards - this lets us set both the run text/formatting
# and the paragraph style (via w:pPr/w:pStyle, or its absence) in one
# shot. Paragraph indices/count are stable across InsertXML calls since
# it replaces a range's contents without adding or removing paragraphs.

$xml6 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Vetores. Vetores no R2  e no R3. Dependência Linear. Produtos de Vetores. A Reta. O Plano. Distâncias. Coordenadas Polares. Mudança de Coordenadas. Cônicas. Superfícies Quádricas. Equações Paramétricas.</w:t></w:r></w:p>
'@
$d.Paragraphs(6).Range.InsertXML($xml6)

$xml7 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:i/></w:rPr><w:t>Vectors. Vectors in 2 and 3 Dimensions. Linear Dependence. Products of Vectors. Lines. Planes. Distances. Polar Coordinates. Coordinates changing. Conic Sections. Quadric Surfaces.</w:t></w:r></w:p>
'@
$d.Paragraphs(7).Range.InsertXML($xml7)

$xml9 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Fornecer fundamentos teóricos sobre vetores, retas no espaço e plano (com suas relações), cônicas e quádricas, tópicos essenciais no estudo de todas Engenharias</w:t></w:r></w:p>
'@
$d.Paragraphs(9).Range.InsertXML($xml9)

$xml11 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>•Vetores: Reta orientada. Eixo. Segmento orientado. Segmentos equipolentes.  Vetor. Operações com vetores. Ângulo de dois vetores.</w:t><w:br/><w:t>•Vetores no r2 e no r3: Decomposição de um vetor no plano. Expressão analítica de um vetor. Igualdade e operações; Vetor definido pelas coordenadas da origem e da extremidade. Decomposição de um vetor no espaço. Igualdade – Operações – Vetor definido pelos pontos extremos. Condição de paralelismo de dois vetores.</w:t><w:br/><w:t>•Dependência linear: Dependência e Independência Linear de vetores no R2 e no R3. Base. Mudança de Base.</w:t><w:br/><w:t>•Produtos de vetores: Produto escalar. Módulo de um vetor. Propriedades do produto escalar. Ângulo de dois vetores. Ângulos diretores e cosenos diretores de um vetor. Projeção de um vetor. Produto escalar no R2. Produto vetorial. Propriedades do produto vetorial. Interpretação geométrica do módulo do produto vetorial de dois vetores. Produto misto. Propriedades do produto misto. Interpretação geométrica do módulo do produto misto.</w:t><w:br/><w:t>•A reta: Equação vetorial da reta. Reta definida por dois pontos. Equações paramétricas da reta. Equações simétricas da reta. Equações reduzidas da reta. Retas paralelas aos planos e aos eixos coordenados. Ângulo de duas retas. Condição de paralelismo e de ortogonalidade de duas retas. Condição de coplanaridade de duas retas. Posições relativas de duas retas. Reta ortogonal a duas retas. Ponto que divide um segmento de reta em uma razão dada.</w:t><w:br/><w:t>•O plano: Equação geral do plano. Determinação de um plano. Planos paralelos aos eixos e aos planos coordenados – Casos particulares. Equações paramétricas do plano. Ângulo de dois planos. Ângulo de uma reta com um plano. Intersecção de dois planos. Intersecção de reta com plano.</w:t><w:br/><w:t>•Distâncias: Distância entre dois pontos. Distância de um ponto a uma reta. Distância de duas retas. Distância de um ponto a um plano. Distância entre dois planos. Distância de uma reta a um plano.</w:t><w:br/><w:t>•Coordenadas polares: Definição de Coordenadas polares, equações e gráficos polares. Relacionando coordenadas polares e coordenadas cartesianas</w:t><w:br/><w:t>•Mudança de coordenadas: Mudança de coordenadas em R2 e em R3. Aplicação de translações e rotações.</w:t><w:br/><w:t>•Equações paramétricas: da reta, da circunferência. Equações Paramétricas de curvas.</w:t><w:br/><w:t>•Cônicas: A parábola. A elipse. A hipérbole. As seções cônicas.</w:t><w:br/><w:t>•Superfícies quádricas: Introdução. Superfícies quádricas centradas. Superfícies quádricas não</w:t></w:r></w:p>
'@
$d.Paragraphs(11).Range.InsertXML($xml11)

$xml12 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:i/></w:rPr><w:t>The discipline aims at providing theoretical foundation regarding vectors, lines and planes (with their relations), conic sections and quadric surfaces, which are essential matters in the study of Engineering.</w:t></w:r></w:p>
'@
$d.Paragraphs(12).Range.InsertXML($xml12)

$xml14 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.</w:t></w:r></w:p>
'@
$d.Paragraphs(14).Range.InsertXML($xml14)

$xml17 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Método: </w:t></w:r><w:r><w:t>NF≥ 5,0.</w:t><w:br/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Critério: </w:t></w:r><w:r><w:t>(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.</w:t><w:br/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Norma de recuperação: </w:t></w:r><w:r><w:t>1.CAMARGO, Ivan ; BOULOS, Paulo. Geometria Analítica: um tratamento vetorial. São Paulo: Prentice Hall, 2005.</w:t><w:br/><w:t>2.LIMA, Elon Lages de. Geometria analítica e algebra Linear. Rio de Janeiro: SBM SociedadeBrasileira de Matemática,2001. Coleção Matemática Universitária.</w:t><w:br/><w:t>3.CAROLI, Alésio de; CALLIOLI, A.; FEITOSA, Miguel O. Matrizes vetores geometria analítica. São Paulo: Nobel, 1998.</w:t><w:br/><w:t>4.SANTOS, Nathan Moreira dos. Vetores e matrizes: uma introdução à álgebra linear. São Paulo: Thomson, 2007.</w:t></w:r></w:p>
'@
$d.Paragraphs(17).Range.InsertXML($xml17)

$xml19 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>3682251 - Gabrielle Weber Martins</w:t></w:r></w:p>
'@
$d.Paragraphs(19).Range.InsertXML($xml19)
